# Weekly update: insert a new week's worth of price data (2 rows) for
# "Pepino ensalada" / Agricola del Norte S.A. de Arica, pushing the
# existing rows 242:245 down to 244:247, then populate the two newly
# inserted rows (242:243) with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above row 242 (shifts old 242:245 -> 244:247,
# inheriting formatting/number-format from the surrounding rows, same as
# Excel's native Rows.Insert behaviour).
$ws.Rows("242:243").Insert()

# New row 242 - "Primera" quality
$ws.Range("A242").Value = 1
$ws.Range("B242").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C242").Value = "Arica y Parinacota"
$ws.Range("D242").Value = 44595
$ws.Range("E242").Value = 15
$ws.Range("F242").Value = 100112043
$ws.Range("G242").Value = "Pepino ensalada"
$ws.Range("H242").Value = "Sin especificar"
$ws.Range("I242").Value = "Primera"
$ws.Range("J242").Value = 160
$ws.Range("K242").Value = 9000
$ws.Range("L242").Value = 10000
$ws.Range("M242").Value = 9500
$ws.Range("N242").Value = "$/caja 70 unidades"
$ws.Range("O242").Value = "Región de Arica y Parinacota"
$ws.Range("P242").Value = 136
$ws.Range("Q242").Value = 70
$ws.Range("R242").Value = "Hortaliza"

# New row 243 - "Segunda" quality
$ws.Range("A243").Value = 1
$ws.Range("B243").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C243").Value = "Arica y Parinacota"
$ws.Range("D243").Value = 44595
$ws.Range("E243").Value = 15
$ws.Range("F243").Value = 100112043
$ws.Range("G243").Value = "Pepino ensalada"
$ws.Range("H243").Value = "Sin especificar"
$ws.Range("I243").Value = "Segunda"
$ws.Range("J243").Value = 160
$ws.Range("K243").Value = 7000
$ws.Range("L243").Value = 8000
$ws.Range("M243").Value = 7500
$ws.Range("N243").Value = "$/caja 100 unidades"
$ws.Range("O243").Value = "Región de Arica y Parinacota"
$ws.Range("P243").Value = 75
$ws.Range("Q243").Value = 100
$ws.Range("R243").Value = "Hortaliza"
